$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 7 and 8: additional testing tasks
$ws.Range("A7").Value = "Test: Generate a YAML file"
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = "Test: Generate the Ansible command"
$ws.Range("B8").Value = 1

# Reword the description text in E6 (new wording for the meetings line)
$ws.Range("E6").Value = "Time spent during oral meetings, discussion, emails and so on."

# Add the pricing formulas (Cost column): F3 stands alone, F4:F19 share one formula
$ws.Range("F3").Formula = "=B3*30"
$ws.Range("F4:F19").Formula = "=B4*30"

# Update selection to match the saved view state
$ws.Range("E7").Select()
